$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 1.029332457036088
$ws.Range("B11").Value = 0.24128467758209
$ws.Range("B12").Value = -1.459370673031765
$ws.Range("B14").Value = 0.7977233685636995
$ws.Range("B15").Value = 1.027544699529146
$ws.Range("B16").Value = 0.1374695422775655
$ws.Range("B17").Value = -0.9225503716806988
$ws.Range("B18").Value = 0.5557457034887239
$ws.Range("B19").Value = 1.11054283609348
$ws.Range("B20").Value = 0.7305722247131936
$ws.Range("B22").Value = -0.4776197014916527
$ws.Range("B23").Value = 0.04982731217580827
$ws.Range("B24").Value = 0.2651053283564908
